{"js": "// Apply the requirement-document edits described by the commit:\n// \"Cambios en requisitos por inconsistencias\"\n\nconst body = context.document.body;\n\n// ------------------------------------------------------------------\n// 1) \"Las empresas pueden crear chollos...\" bullet:\n//    \" la fecha en la que se cre\u00f3,\" -> \" informaci\u00f3n sobre si est\u00e1 publicado,\"\n// ------------------------------------------------------------------\n{\n  const results = body.search(\"la fecha en la que se cre\u00f3,\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"informaci\u00f3n sobre si est\u00e1 publicado,\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// ------------------------------------------------------------------\n// 2) \"El sistema guardar\u00e1 los siguientes par\u00e1metros...\" bullet:\n//    append new clause about default avatar/image after \"comisi\u00f3n por ventas.\"\n// ------------------------------------------------------------------\n{\n  const results = body.search(\"comisi\u00f3n por ventas.\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"comisi\u00f3n por ventas, avatar por defecto e imagen por defecto para URLs de im\u00e1genes que no existan.\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// ------------------------------------------------------------------\n// 3) Delete the stand-alone bullet:\n//    \"Eliminar etiquetas que considere inapropiadas o que est\u00e9n en desuso.\"\n// ------------------------------------------------------------------\n{\n  const results = body.search(\n    \"Eliminar etiquetas que considere inapropiadas o que est\u00e9n en desuso.\",\n    { matchCase: true }\n  );\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    const para = results.items[0].paragraphs.getFirst();\n    para.delete();\n    await context.sync();\n  }\n}\n\n// ------------------------------------------------------------------\n// 4) \"Un chollo no ser\u00e1 p\u00fablico hasta dos d\u00edas despu\u00e9s de su creaci\u00f3n,\n//    solo visible durante dichos dos d\u00edas para los usuarios \" ->\n//    \"Un chollo no p\u00fablico solo ser\u00e1 visible para los usuarios \"\n//    and relocate the \"_GoBack\" bookmark into this paragraph, right\n//    after \"visible \".\n// ------------------------------------------------------------------\n{\n  const results = body.search(\n    \"Un chollo no ser\u00e1 p\u00fablico hasta dos d\u00edas despu\u00e9s de su creaci\u00f3n, solo visible durante dichos dos d\u00edas para los usuarios \",\n    { matchCase: true }\n  );\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"Un chollo no p\u00fablico solo ser\u00e1 visible para los usuarios \",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n\n  // The document only ever has a single \"_GoBack\" bookmark; remove it\n  // from its old location before re-adding it in the new spot.\n  const oldBookmark = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\n  oldBookmark.load(\"isNullObject\");\n  await context.sync();\n  if (!oldBookmark.isNullObject) {\n    context.document.deleteBookmark(\"_GoBack\");\n    await context.sync();\n  }\n\n  const newSpot = body.search(\"Un chollo no p\u00fablico solo ser\u00e1 visible \", {\n    matchCase: true,\n  });\n  newSpot.load(\"items\");\n  await context.sync();\n  if (newSpot.items.length > 0) {\n    const collapsed = newSpot.items[0].getRange(\"End\");\n    collapsed.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# Apply the requirement-document edits described by the commit:\n# \"Cambios en requisitos por inconsistencias\"\n\n$d = $word.ActiveDocument\n\n# ------------------------------------------------------------------\n# 1) \"Las empresas pueden crear chollos...\" bullet:\n#    \" la fecha en la que se cre\u00f3,\" -> \" informaci\u00f3n sobre si est\u00e1 publicado,\"\n# ------------------------------------------------------------------\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Replacement.ClearFormatting()\n$rng1.Find.Execute(\n    \"la fecha en la que se cre\u00f3,\",  # FindText\n    $false,                          # MatchCase\n    $false,                          # MatchWholeWord\n    $false,                          # MatchWildcards\n    $false,                          # MatchSoundsLike\n    $false,                          # MatchAllWordForms\n    $true,                           # Forward\n    1,                                # Wrap (wdFindContinue)\n    $false,                          # Format\n    \"informaci\u00f3n sobre si est\u00e1 publicado,\",  # ReplaceWith\n    2                                 # Replace (wdReplaceAll)\n) | Out-Null\n\n# ------------------------------------------------------------------\n# 2) \"El sistema guardar\u00e1 los siguientes par\u00e1metros...\" bullet:\n#    append new clause about default avatar/image after \"comisi\u00f3n por ventas.\"\n# ------------------------------------------------------------------\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Replacement.ClearFormatting()\n$rng2.Find.Execute(\n    \"comisi\u00f3n por ventas.\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"comisi\u00f3n por ventas, avatar por defecto e imagen por defecto para URLs de im\u00e1genes que no existan.\",\n    2\n) | Out-Null\n\n# ------------------------------------------------------------------\n# 3) Delete the stand-alone bullet:\n#    \"Eliminar etiquetas que considere inapropiadas o que est\u00e9n en desuso.\"\n# ------------------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"Eliminar etiquetas que considere inapropiadas o que est\u00e9n en desuso.\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# ------------------------------------------------------------------\n# 4) \"Un chollo no ser\u00e1 p\u00fablico hasta dos d\u00edas despu\u00e9s de su creaci\u00f3n,\n#    solo visible durante dichos dos d\u00edas para los usuarios \" ->\n#    \"Un chollo no p\u00fablico solo ser\u00e1 visible para los usuarios \"\n#    and relocate the \"_GoBack\" bookmark into this paragraph, right\n#    after \"visible \".\n# ------------------------------------------------------------------\n$rng4 = $d.Content\n$rng4.Find.ClearFormatting()\n$rng4.Find.Replacement.ClearFormatting()\n$rng4.Find.Execute(\n    \"Un chollo no ser\u00e1 p\u00fablico hasta dos d\u00edas despu\u00e9s de su creaci\u00f3n, solo visible durante dichos dos d\u00edas para los usuarios \",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Un chollo no p\u00fablico solo ser\u00e1 visible para los usuarios \",\n    2\n) | Out-Null\n\n# The document only ever has a single \"_GoBack\" bookmark; remove it from\n# its old location before re-adding it in the new spot.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$rngBm = $d.Content\n$rngBm.Find.ClearFormatting()\n$rngBm.Find.Text = \"Un chollo no p\u00fablico solo ser\u00e1 visible \"\n$rngBm.Find.Execute() | Out-Null\n$rngBm.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $rngBm) | Out-Null\n"}
